# The source diff for this workbook (Data/EC/NIT-9000075545.xlsx) only
# touches bookkeeping that Excel regenerates on every save:
#   - xl/workbook.xml: fileVersion/rupBuild, xr:revisionPtr documentId,
#     bookViews xr2:uid  -> random session/revision GUIDs + build number.
#   - xl/sharedStrings.xml: the ten period labels ("1608".."1705") are
#     physically reordered inside <sst>, but every <c> in sheet1.xml keeps
#     referencing the very same <si> index it did before, so the values
#     shown in B15:J25 do not change at all.
#   - xl/styles.xml: borders 1 and 2 swap their <right> side, and every
#     cellXfs entry that pointed at borderId 1/2 swaps to the other index
#     in lock-step, so every cell keeps the exact same rendered border.
#   - xl/drawings/drawing1.xml: the picture's a16:creationId GUID changes.
#   - xl/worksheets/sheet1.xml: only the xr:uid attribute (a random
#     per-sheet revision GUID) changes; <sheetData> is byte-for-byte
#     identical before and after.
#
# None of that is reachable (or meaningful) through the Excel object
# model - it is exactly what Excel rewrites whenever a workbook is
# opened and re-saved, here under a newer Office build (rupBuild
# 29127 -> 29231), without the user actually changing any cell value,
# formula, or visible formatting. So the correct replay is simply to
# touch the workbook/sheet and leave its content untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
